$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look numeric stay exactly as text
# (preserve trailing zeros / exact formatting as in source data).
$ws.Range('D2').Value = '26.924.92'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '1.817.65'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4650'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3706'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07351'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8711'
$ws.Range('E10').Value = '  +0.13%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = '1.829.91'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07089'
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.514'
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008718'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.72'
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('D21').Value = '26.943.03'
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.326'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('E23').Value = '  -3.21%  '
$ws.Range('D24').Value = '2.054.35'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('E25').Value = '  -2.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.02'
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.43'
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.139'
$ws.Range('E28').Value = '  -4.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.311'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.56'
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08907'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7595'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.154'
$ws.Range('E33').Value = '  -2.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.465'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.921'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.096'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01958'
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05263'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.937'
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.246'
$ws.Range('E41').Value = '  -0.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5337'
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('E44').Value = '  -1.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.458'
$ws.Range('E45').Value = '  -1.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4943'
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.33'
$ws.Range('E47').Value = '  -0.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.680'
$ws.Range('E48').Value = '  +0.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.001'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.33'
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06287'
$ws.Range('E51').Value = '  -0.64%  '
